$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.022.02'
$ws.Range("E2").Value = '  -1.91%  '

$ws.Range("D3").Value = '2.466.72'
$ws.Range("E3").Value = '  -2.38%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = "'518.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.43%  '

$ws.Range("D6").Value = "'131.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.40%  '

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -2.19%  '

$ws.Range("D9").Value = "'0.0993"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.03%  '

$ws.Range("D11").Value = "'5.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.12%  '

$ws.Range("E12").Value = '  -1.65%  '

$ws.Range("D13").Value = '2.900.01'
$ws.Range("E13").Value = '  -2.51%  '

$ws.Range("D14").Value = '57.874.62'
$ws.Range("E14").Value = '  -2.05%  '

$ws.Range("D15").Value = "'22.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.26%  '

$ws.Range("E16").Value = '  -2.49%  '

$ws.Range("D17").Value = '2.459.78'
$ws.Range("E17").Value = '  -2.70%  '

$ws.Range("D18").Value = "'10.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.48%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = "'4.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.48%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = "'319.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.17%  '

$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").Value = "'5.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.58%  '

$ws.Range("D23").Value = "'63.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.79%  '

$ws.Range("D24").Value = "'0.411"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.78%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D27").Value = "'7.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.05%  '

$ws.Range("D28").Value = '0.0₃0752'
$ws.Range("E28").Value = '  -3.03%  '

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = "'166.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.75%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'1.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.59%  '

$ws.Range("D31").Value = "'6.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.08%  '

$ws.Range("E32").Value = '  -2.12%  '

$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("D34").Value = "'0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.13%  '

$ws.Range("D35").Value = "'18.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.20%  '

$ws.Range("D36").Value = "'1.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -10.63%  '

$ws.Range("D37").Value = "'3.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.20%  '

$ws.Range("D38").Value = "'1.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.58%  '

$ws.Range("D39").Value = "'0.790"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.26%  '

$ws.Range("D40").Value = "'3.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.57%  '

$ws.Range("D41").Value = "'272.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.45%  '

$ws.Range("D42").Value = "'5.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.67%  '

$ws.Range("E43").Value = '  -2.73%  '

$ws.Range("D44").Value = "'126.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.57%  '

$ws.Range("D45").Value = "'0.0905"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.21%  '

$ws.Range("E46").Value = '  -4.34%  '

$ws.Range("E47").Value = '  -3.53%  '

$ws.Range("D48").Value = "'17.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.07%  '

$ws.Range("D49").Value = '1.730.16'
$ws.Range("E49").Value = '  -2.03%  '

$ws.Range("E50").Value = '  -1.52%  '

$ws.Range("D51").Value = "'4.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.18%  '
